$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: After the first paragraph ("Play Colossal Gems Slot Free - Low
# Volatility and Giant Gem Symbols", Heading 1), insert a new plain paragraph
# containing a "Meta description" label (bold) followed by the description
# text (normal weight), mirroring the leading empty-run pattern used
# elsewhere in the document.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Colossal Gems slot game with low volatility Respin and Free Spin features with giant gem symbols. Play for free on any device.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# Change 2: Near the end of the document, remove the duplicated bold title
# paragraph entirely, and replace the text of the following italic paragraph
# (originally the meta-description recap) with the new DALLE image prompt,
# keeping its italic formatting and leading empty run intact.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePara = $null
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Play Colossal Gems Slot Free - Low Volatility and Giant Gem Symbols") {
        $dupTitlePara = $p
        break
    }
}

$dupTitlePara.Range.Delete()

$count2 = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($count2)
$promptStart = $promptPara.Range.Start
$promptEnd = $promptPara.Range.End
$promptRange = $d.Range($promptStart, $promptEnd)

$promptRange.Text = 'Prompt: Create a feature image for "Colossal Gems" in a cartoon style featuring a happy Maya warrior with glasses. DALLE, create a vibrant and colorful feature image for "Colossal Gems" that captures the essence of this joyful gem-themed slot game. Make sure to include a happy Maya warrior with glasses in the image to highlight its adventurous and playful nature. Use bright colors like pink, orange, and purple to give the image an eye-catching and dynamic look. Incorporate giant gem symbols of mega sizes, such as 2x2, 3x3, and 4x4, in the background to showcase the exciting bonus features of the game. Make it fun and inviting to encourage players to try their luck with this low volatility slot.'
